# Updates the cryptos list (prices/24h volume%) to the latest scraped values.
# Mirrors commit: "Updated cryptos list on Sun Aug 20 20:47:52 UTC 2023 with GitHub Actions"
# Rows 31/32 (Filecoin / InternetComputer(DFINITY)) also swapped rank order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# List of cell updates: Cell address -> new text value.
$updates = @(
    @{ Cell = 'D2'; Value = '26.482.71' }
    @{ Cell = 'E2'; Value = '  +0.67%  ' }
    @{ Cell = 'D3'; Value = '1.704.62' }
    @{ Cell = 'E3'; Value = '  +1.19%  ' }
    @{ Cell = 'D4'; Value = '1.009' }
    @{ Cell = 'E4'; Value = '  +0.04%  ' }
    @{ Cell = 'D5'; Value = '219.55' }
    @{ Cell = 'E5'; Value = '  +0.70%  ' }
    @{ Cell = 'D6'; Value = '0.5508' }
    @{ Cell = 'E6'; Value = '  +5.37%  ' }
    @{ Cell = 'D7'; Value = '1.009' }
    @{ Cell = 'E7'; Value = '  -0.06%  ' }
    @{ Cell = 'D8'; Value = '0.2747' }
    @{ Cell = 'E8'; Value = '  +1.62%  ' }
    @{ Cell = 'D9'; Value = '0.06489' }
    @{ Cell = 'E9'; Value = '  +1.38%  ' }
    @{ Cell = 'D10'; Value = '22.12' }
    @{ Cell = 'E10'; Value = '  +0.64%  ' }
    @{ Cell = 'D11'; Value = '0.07708' }
    @{ Cell = 'E11'; Value = '  +2.73%  ' }
    @{ Cell = 'D12'; Value = '4.560' }
    @{ Cell = 'E12'; Value = '  -0.06%  ' }
    @{ Cell = 'D13'; Value = '1.687.28' }
    @{ Cell = 'E13'; Value = '  -1.22%  ' }
    @{ Cell = 'D14'; Value = '0.5858' }
    @{ Cell = 'E14'; Value = '  +1.24%  ' }
    @{ Cell = 'D15'; Value = '0.000008420' }
    @{ Cell = 'E15'; Value = '  -0.06%  ' }
    @{ Cell = 'D16'; Value = '66.02' }
    @{ Cell = 'E16'; Value = '  +2.77%  ' }
    @{ Cell = 'D17'; Value = '26.518.59' }
    @{ Cell = 'E17'; Value = '  +0.56%  ' }
    @{ Cell = 'D18'; Value = '4.968' }
    @{ Cell = 'E18'; Value = '  +1.10%  ' }
    @{ Cell = 'D19'; Value = '1.009' }
    @{ Cell = 'E19'; Value = '  +0.07%  ' }
    @{ Cell = 'D20'; Value = '11.00' }
    @{ Cell = 'E20'; Value = '  +1.35%  ' }
    @{ Cell = 'D21'; Value = '193.01' }
    @{ Cell = 'E21'; Value = '  +2.43%  ' }
    @{ Cell = 'D22'; Value = '6.283' }
    @{ Cell = 'E22'; Value = '  +1.65%  ' }
    @{ Cell = 'D23'; Value = '1.009' }
    @{ Cell = 'E23'; Value = '  -0.04%  ' }
    @{ Cell = 'D24'; Value = '149.05' }
    @{ Cell = 'E24'; Value = '  +3.11%  ' }
    @{ Cell = 'D25'; Value = '0.1335' }
    @{ Cell = 'E25'; Value = '  +8.66%  ' }
    @{ Cell = 'D26'; Value = '7.944' }
    @{ Cell = 'E26'; Value = '  +3.35%  ' }
    @{ Cell = 'D27'; Value = '15.92' }
    @{ Cell = 'E27'; Value = '  +0.81%  ' }
    @{ Cell = 'D28'; Value = '0.06317' }
    @{ Cell = 'E28'; Value = '  -4.57%  ' }
    @{ Cell = 'D29'; Value = '1.382' }
    @{ Cell = 'E29'; Value = '  +2.82%  ' }
    @{ Cell = 'E30'; Value = '  +0.35%  ' }
    @{ Cell = 'B31'; Value = 'InternetComputer(DFINITY)' }
    @{ Cell = 'C31'; Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp' }
    @{ Cell = 'D31'; Value = '3.626' }
    @{ Cell = 'E31'; Value = '  +1.74%  ' }
    @{ Cell = 'B32'; Value = 'Filecoin' }
    @{ Cell = 'C32'; Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil' }
    @{ Cell = 'D32'; Value = '3.615' }
    @{ Cell = 'E32'; Value = '  +1.31%  ' }
    @{ Cell = 'D33'; Value = '1.697' }
    @{ Cell = 'E33'; Value = '  +2.50%  ' }
    @{ Cell = 'D34'; Value = '1.048' }
    @{ Cell = 'E34'; Value = '  +2.25%  ' }
    @{ Cell = 'D35'; Value = '0.6209' }
    @{ Cell = 'E35'; Value = '  +0.45%  ' }
    @{ Cell = 'D36'; Value = '2.409' }
    @{ Cell = 'E36'; Value = '  +0.35%  ' }
    @{ Cell = 'E37'; Value = '  +2.58%  ' }
    @{ Cell = 'D38'; Value = '0.01652' }
    @{ Cell = 'E38'; Value = '  +2.45%  ' }
    @{ Cell = 'D39'; Value = '1.122.67' }
    @{ Cell = 'E39'; Value = '  +1.47%  ' }
    @{ Cell = 'D40'; Value = '6.174' }
    @{ Cell = 'E40'; Value = '  -3.15%  ' }
    @{ Cell = 'D41'; Value = '0.8862' }
    @{ Cell = 'E41'; Value = '  +0.99%  ' }
    @{ Cell = 'D42'; Value = '1.017' }
    @{ Cell = 'E42'; Value = '  +0.05%  ' }
    @{ Cell = 'D43'; Value = '101.57' }
    @{ Cell = 'E43'; Value = '  +0.75%  ' }
    @{ Cell = 'E44'; Value = '  +1.05%  ' }
    @{ Cell = 'D45'; Value = '57.95' }
    @{ Cell = 'E45'; Value = '  +2.31%  ' }
    @{ Cell = 'D46'; Value = '0.00000000109' }
    @{ Cell = 'E46'; Value = '  -2.01%  ' }
    @{ Cell = 'D47'; Value = '8.240' }
    @{ Cell = 'E47'; Value = '  +0.76%  ' }
    @{ Cell = 'E48'; Value = '  -0.17%  ' }
    @{ Cell = 'D49'; Value = '0.05283' }
    @{ Cell = 'E49'; Value = '  +0.24%  ' }
    @{ Cell = 'D50'; Value = '6.155' }
    @{ Cell = 'E50'; Value = '  +1.96%  ' }
    @{ Cell = 'D51'; Value = '0.4303' }
    @{ Cell = 'E51'; Value = '  -0.16%  ' }
)

foreach ($u in $updates) {
    $range = $ws.Range($u.Cell)
    $text = $u.Value
    # Cells hold text values (e.g. prices like "1.009" or "26.425.89").
    # Prefix ambiguous numeric-looking text with a quote so Excel keeps it as
    # a string instead of silently converting it to a Number, then restore the
    # default "Normal" style so no extra number formatting is left behind.
    $needsQuote = $false
    if ($text -match "^[+-]?[0-9]*\.?[0-9]+([eE][+-]?[0-9]+)?$") {
        $needsQuote = $true
    }
    if ($needsQuote) {
        $range.Value = "'" + $text
        $range.Style = "Normal"
    } else {
        $range.Value = $text
    }
}
